# Qatar Stars League workbook update
# The "Atualização de bases das ligas" commit swaps the full data payload
# (columns B:AD) between pairs of rows, while leaving the row-index
# column A untouched. Swap each pair using Range.Value2 (Value is not
# usable through this COM shim for ranges/cells - it resolves to the
# indexed property descriptor instead of the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$pairs = @(
    @(21, 22),
    @(24, 25),
    @(40, 41),
    @(42, 43),
    @(45, 46),
    @(50, 51),
    @(68, 69),
    @(81, 82),
    @(83, 84),
    @(87, 88),
    @(94, 95),
    @(102, 103),
    @(110, 111),
    @(112, 113),
    @(118, 119),
    @(122, 124),
    @(125, 127),
    @(129, 130)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$($r1):AD$($r1)")
    $range2 = $ws.Range("B$($r2):AD$($r2)")

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value2 = $v2
    $range2.Value2 = $v1
}
